$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Colname 1","Colname 2","Colname 1","Colname 2","Colname 1","Colname 2","Colname 1","Colname 2","Colname 3","Colname 4","Colname 5","Colname 6"),
    @("cell1","cell2","cell3","cell4","cell5","cell6","cell7","cell8","cell9","cell10","cell11","cell5"),
    @("cell3","cell4","cell3","cell4","cell3","cell4","cell3","cell4","cell3","cell4","cell3","cell4"),
    @("cell5","cell6","cell5","cell6","cell5","cell6","cell5","cell6","cell5","cell6","cell5","cell6"),
    @("cell7","cell8","cell7","cell8","cell7","cell8","cell7","cell8","cell7","cell8","cell7","cell8"),
    @("cell9","cell10","cell9","cell10","cell9","cell10","cell9","cell10","cell9","cell10","cell9","cell10"),
    @("cell11","cell12","cell11","cell12","cell11","cell12","cell11","cell12","cell11","cell12","cell11","cell12"),
    @("cell13","cell14","cell13","cell14","cell13","cell14","cell13","cell14","cell13","cell14","cell13","cell14")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Range("L2").Select()
